# Swap the presentation's applied theme color scheme from the
# "Integral" (Red Violet) palette over to the stock "Office Theme"
# (Office) palette.
#
# ppt/theme/theme1.xml (the theme used by the slide master / every
# slide) originally carried the "Integral" / "Red Violet" color
# scheme; the edit turns it into the standard "Office" palette used
# by ppt/theme/theme2.xml (the notes-master theme). The PowerPoint
# object model exposes the live theme's 12 scheme colors through
# Theme.ThemeColorScheme(Index).RGB, in the fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
#   11 hlink, 12 folHlink
# RGB values use the standard COM encoding R + G*256 + B*65536.

function ConvertTo-ComRgb($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# Target "Office Theme" / "Office" scheme colors, in
# ThemeColorScheme.Item index order.
$officeThemeColors = @(
    (ConvertTo-ComRgb 0x00 0x00 0x00), # 1  dk1       000000
    (ConvertTo-ComRgb 0xFF 0xFF 0xFF), # 2  lt1       FFFFFF
    (ConvertTo-ComRgb 0x44 0x54 0x6A), # 3  dk2       44546A
    (ConvertTo-ComRgb 0xE7 0xE6 0xE6), # 4  lt2       E7E6E6
    (ConvertTo-ComRgb 0x5B 0x9B 0xD5), # 5  accent1   5B9BD5
    (ConvertTo-ComRgb 0xED 0x7D 0x31), # 6  accent2   ED7D31
    (ConvertTo-ComRgb 0xA5 0xA5 0xA5), # 7  accent3   A5A5A5
    (ConvertTo-ComRgb 0xFF 0xC0 0x00), # 8  accent4   FFC000
    (ConvertTo-ComRgb 0x44 0x72 0xC4), # 9  accent5   4472C4
    (ConvertTo-ComRgb 0x70 0xAD 0x47), # 10 accent6   70AD47
    (ConvertTo-ComRgb 0x05 0x63 0xC1), # 11 hlink     0563C1
    (ConvertTo-ComRgb 0x95 0x4F 0x72)  # 12 folHlink  954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
